$wb = $excel.ActiveWorkbook

# --- Users sheet: replace the sample user row (Ashley Choi -> Amy Xia) ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Amy Xia"
$wsUsers.Range("B2").Value = "Time Tracking Litigation "
$wsUsers.Range("C2").Value = "CF Financial User"

# --- TitleRateSheet: update selection, then leave it inactive ---
$wsTitleRate = $wb.Worksheets.Item("TitleRateSheet")
$wsTitleRate.Activate() | Out-Null
$wsTitleRate.Range("J22").Select() | Out-Null

# --- Users sheet becomes the active tab with a new selection ---
$wsUsers.Activate() | Out-Null
$wsUsers.Range("C17").Select() | Out-Null
